$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This weekly update prepends 4 new price-report rows (dated 45077) ahead of
# the existing chronological block that starts at row 880, pushing the old
# rows 880-990 down to 884-994 (dimension grows from A1:R990 to A1:R994).
$ws.Range("A880:R883").EntireRow.Insert()

# Build the 4 new rows (columns A-R) as a single 2D array and write them in
# one shot into the freshly inserted, currently-blank rows.
$newRows = New-Object 'object[,]' 4,18

# Row 880: Primera, Región Metropolitana
$newRows[0,0]  = 9
$newRows[0,1]  = "Vega Central Mapocho de Santiago"
$newRows[0,2]  = "Metropolitana"
$newRows[0,3]  = 45077
$newRows[0,4]  = 13
$newRows[0,5]  = 100112023
$newRows[0,6]  = "Brócoli"
$newRows[0,7]  = "Sin especificar"
$newRows[0,8]  = "Primera"
$newRows[0,9]  = 3400
$newRows[0,10] = 800
$newRows[0,11] = 900
$newRows[0,12] = 850
$newRows[0,13] = "`$/unidad"
$newRows[0,14] = "Región Metropolitana"
$newRows[0,15] = 850
$newRows[0,16] = 1
$newRows[0,17] = "Hortaliza"

# Row 881: Primera, Región de O'Higgins
$newRows[1,0]  = 9
$newRows[1,1]  = "Vega Central Mapocho de Santiago"
$newRows[1,2]  = "Metropolitana"
$newRows[1,3]  = 45077
$newRows[1,4]  = 13
$newRows[1,5]  = 100112023
$newRows[1,6]  = "Brócoli"
$newRows[1,7]  = "Sin especificar"
$newRows[1,8]  = "Primera"
$newRows[1,9]  = 1600
$newRows[1,10] = 800
$newRows[1,11] = 900
$newRows[1,12] = 850
$newRows[1,13] = "`$/unidad"
$newRows[1,14] = "Región de O'Higgins"
$newRows[1,15] = 850
$newRows[1,16] = 1
$newRows[1,17] = "Hortaliza"

# Row 882: Segunda, Región Metropolitana
$newRows[2,0]  = 9
$newRows[2,1]  = "Vega Central Mapocho de Santiago"
$newRows[2,2]  = "Metropolitana"
$newRows[2,3]  = 45077
$newRows[2,4]  = 13
$newRows[2,5]  = 100112023
$newRows[2,6]  = "Brócoli"
$newRows[2,7]  = "Sin especificar"
$newRows[2,8]  = "Segunda"
$newRows[2,9]  = 1690
$newRows[2,10] = 700
$newRows[2,11] = 700
$newRows[2,12] = 700
$newRows[2,13] = "`$/unidad"
$newRows[2,14] = "Región Metropolitana"
$newRows[2,15] = 700
$newRows[2,16] = 1
$newRows[2,17] = "Hortaliza"

# Row 883: Segunda, Región de O'Higgins
$newRows[3,0]  = 9
$newRows[3,1]  = "Vega Central Mapocho de Santiago"
$newRows[3,2]  = "Metropolitana"
$newRows[3,3]  = 45077
$newRows[3,4]  = 13
$newRows[3,5]  = 100112023
$newRows[3,6]  = "Brócoli"
$newRows[3,7]  = "Sin especificar"
$newRows[3,8]  = "Segunda"
$newRows[3,9]  = 970
$newRows[3,10] = 700
$newRows[3,11] = 700
$newRows[3,12] = 700
$newRows[3,13] = "`$/unidad"
$newRows[3,14] = "Región de O'Higgins"
$newRows[3,15] = 700
$newRows[3,16] = 1
$newRows[3,17] = "Hortaliza"

$ws.Range("A880:R883").Value2 = $newRows

# Column D (Fecha) carries a date number format (style index 2 in the
# original file); make sure the newly written rows keep that same
# date/time format like every other row in this column.
$ws.Range("D880:D883").NumberFormat = $ws.Range("D884").NumberFormat
